$d = $word.ActiveDocument

$pairs = @(
    @("657÷6=109, 3", "236÷7=33, 5"),
    @("798÷6=133, 0", "836÷5=167, 1"),
    @("849÷2=424, 1", "193÷2=96, 1"),
    @("397÷9=44, 1", "176÷8=22, 0"),
    @("366÷6=61, 0", "395÷2=197, 1"),
    @("866÷4=216, 2", "877÷4=219, 1"),
    @("260÷2=130, 0", "543÷7=77, 4"),
    @("748÷4=187, 0", "728÷5=145, 3"),
    @("379÷2=189, 1", "671÷4=167, 3"),
    @("325÷4=81, 1", "113÷6=18, 5"),
    @("671÷6=111, 5", "199÷3=66, 1"),
    @("171÷9=19, 0", "529÷8=66, 1"),
    @("879÷2=439, 1", "477÷7=68, 1"),
    @("506÷7=72, 2", "350÷5=70, 0"),
    @("173÷7=24, 5", "420÷9=46, 6"),
    @("278÷5=55, 3", "192÷9=21, 3"),
    @("873÷5=174, 3", "386÷6=64, 2"),
    @("836÷2=418, 0", "917÷9=101, 8"),
    @("930÷3=310, 0", "281÷5=56, 1"),
    @("725÷8=90, 5", "446÷9=49, 5"),
    @("347÷5=69, 2", "234÷3=78, 0"),
    @("168÷3=56, 0", "531÷3=177, 0"),
    @("790÷9=87, 7", "153÷9=17, 0"),
    @("241÷3=80, 1", "171÷8=21, 3"),
    @("731÷7=104, 3", "361÷5=72, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
